# Add "copy editing" to the skills sheet (soft_skill), bump the order of the
# skills that come after it, and make the skills sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skills")

# The existing soft_skill rows are (in sheet order):
#   24: technical writing               order 1
#   25: cross-functional communication  order 2 -> 3
#   26: stakeholder engagement          order 3 -> 4
# "copy editing" is inserted as the new order-2 soft skill, so bump the two
# that follow it down by one.
$ws.Range("C25").Value2 = 3
$ws.Range("C26").Value2 = 4

# Append the new skill row.
$ws.Range("A27").Value2 = "copy editing"
$ws.Range("B27").Value2 = "soft_skill"
$ws.Range("C27").Value2 = 2

# Make "skills" the active sheet/tab, with the selection sitting just below
# the newly added row.
$ws.Activate()
$ws.Range("B28").Select()
